$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01459712982177734
$ws.Range("C2").Value = 0.02440824508666992
$ws.Range("D2").Value = 0.006189537048339844
$ws.Range("E2").Value = 0.01681041717529297
$ws.Range("F2").Value = 0.002393722534179688
$ws.Range("G2").Value = 0.0571098804473877
$ws.Range("H2").Value = 0.01439552307128906
$ws.Range("I2").Value = 0.02533020973205567
$ws.Range("J2").Value = 0.01119651794433594
$ws.Range("K2").Value = 0.02079572677612305
$ws.Range("L2").Value = 0.003200387954711914
$ws.Range("M2").Value = 0.01460742950439453
$ws.Range("B3").Value = 0.07539291381835937
$ws.Range("C3").Value = 0.02686405181884766
$ws.Range("D3").Value = 0.01516590118408203
$ws.Range("E3").Value = 0.0101994514465332
$ws.Range("F3").Value = 0.01186113357543945
$ws.Range("G3").Value = 0.01146078109741211
$ws.Range("H3").Value = 0.1516443252563477
$ws.Range("I3").Value = 0.04450893402099609
$ws.Range("J3").Value = 0.1031134128570557
$ws.Range("K3").Value = 0.03127560615539551
$ws.Range("L3").Value = 0.03315143585205078
$ws.Range("M3").Value = 0.01675300598144531
$ws.Range("B4").Value = 0.04499831199645996
$ws.Range("C4").Value = 0.02724075317382812
$ws.Range("D4").Value = 0.02290120124816895
$ws.Range("E4").Value = 0.01682343482971192
$ws.Range("F4").Value = 0.0870091438293457
$ws.Range("G4").Value = 0.0116429328918457
$ws.Range("H4").Value = 0.03433222770690918
$ws.Range("I4").Value = 0.02600207328796387
$ws.Range("J4").Value = 0.03202948570251465
$ws.Range("K4").Value = 0.02023811340332031
$ws.Range("L4").Value = 0.05199732780456543
$ws.Range("M4").Value = 0.01557474136352539
$ws.Range("B5").Value = 0.03359456062316894
$ws.Range("C5").Value = 0.02798910140991211
$ws.Range("D5").Value = 0.02648453712463379
$ws.Range("E5").Value = 0.0227330207824707
$ws.Range("H5").Value = 0.02430157661437988
$ws.Range("I5").Value = 0.02726325988769531
$ws.Range("J5").Value = 0.02553739547729492
$ws.Range("K5").Value = 0.02647767066955566
$ws.Range("B6").Value = 0.8559419631958007
$ws.Range("C6").Value = 0.1351036548614502
$ws.Range("D6").Value = 0.6783699989318848
$ws.Range("E6").Value = 0.1171733379364014
$ws.Range("F6").Value = 0.1928281784057617
$ws.Range("G6").Value = 0.06292157173156739
$ws.Range("H6").Value = 1.036152076721191
$ws.Range("I6").Value = 0.1854721069335938
$ws.Range("J6").Value = 0.6953098773956299
$ws.Range("K6").Value = 0.1617080688476563
$ws.Range("L6").Value = 0.2982239246368408
$ws.Range("M6").Value = 0.06996550559997558
